# Fruta / hortaliza, semanal
# Insert the new weekly record as row 5, pushing the existing rows 5-10 down
# to 6-11 (matching the shift seen for every row after the insertion point).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(5).Insert()

$ws.Range("A5").Value = 10
$ws.Range("B5").Value = "Vega Modelo de Temuco"
$ws.Range("C5").Value = "La Araucanía"
$ws.Range("D5").Value = 44601
$ws.Range("E5").Value = 9
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100104
$ws.Range("H5").Value = "Frutos de pepita"
$ws.Range("I5").Value = 100104005
$ws.Range("J5").Value = "Pera asiática"
$ws.Range("K5").Value = "Hosui"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 30
$ws.Range("N5").Value = 28000
$ws.Range("O5").Value = 28000
$ws.Range("P5").Value = 28000
$ws.Range("Q5").Value = "$/caja 18 kilos granel"
$ws.Range("R5").Value = "Región de O'Higgins"
$ws.Range("S5").Value = 1556
$ws.Range("T5").Value = 18
